$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet (was "vestfoldlab_to_vannmiljo")
$ws.Name = "to_vannmiljo"

# Reflect the user's last selection on the sheet before saving
$ws.Range("C29").Select() | Out-Null
